# Supplemental calcs: energy recovery / landfill cost assumptions added to Sheet1 (rows 30-39)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$epaUrl = "https://nepis.epa.gov/Exe/ZyNET.exe/P10053DP.txt?ZyActionD=ZyDocument&Client=EPA&Index=2000%20Thru%202005&Docs=&Query=%28cost%29%20OR%20FNAME%3D%22P10053DP.txt%22%20AND%20FNAME%3D%22P10053DP.txt%22&Time=&EndTime=&SearchMethod=1&TocRestrict=n&Toc=&TocEntry=&QField=&QFieldYear=&QFieldMonth=&QFieldDay=&UseQField=&IntQFieldOp=0&ExtQFieldOp=0&XmlQuery=&File=D%3A%5CZYFILES%5CINDEX%20DATA%5C00THRU05%5CTXT%5C00000021%5CP10053DP.txt&User=ANONYMOUS&Password=anonymous&SortMethod=h%7C-&MaximumDocuments=1&FuzzyDegree=0&ImageQuality=r75g8/r75g8/x150y150g16/i425&Display=hpfr&DefSeekPage=x&SearchBack=ZyActionL&Back=ZyActionS&BackDesc=Results%20page&MaximumPages=1&ZyEntry=7&slide"
$eiaUrl = "https://www.eia.gov/electricity/monthly/epm_table_grapher.php?t=epmt_5_6_a"

# ---- Text/label cells, entered in the same order the original author typed them ----
# (this preserves the shared-string table order of the authored workbook)
$ws.Range("K30").Value = $epaUrl
$ws.Range("K31").Value = $epaUrl
$ws.Range("J30").Value = "`$/Mg wet"
$ws.Range("J31").Value = "% solids"
$ws.Range("J33").Value = "`$/Mg VSS landfilled"
$ws.Range("J34").Value = "`$/kg VSS landfilled"
$ws.Range("J35").Value = "`$/kg VSS Sent to Anaerobic Digester"
$ws.Range("J32").Value = "Fraction of sludge reduced in AD"
$ws.Range("G30").Value = "kgVSS"
$ws.Range("G31").Value = "kgVSS to biogas"
$ws.Range("G32").Value = "kgVSS to landfill"
$ws.Range("G33").Value = "kg wet sludge to landfill"
$ws.Range("J36").Value = "% biogas methane"
$ws.Range("J37").Value = "gVSS/gCOD"
$ws.Range("G34").Value = "cost to landfill 1kg sludge produced "
$ws.Range("G35").Value = "kgCOD to biogas (as methane)"
$ws.Range("J38").Value = "MJ/kg"
$ws.Range("G36").Value = "MJ produced per 1 kg of sludge produced"
$ws.Range("G37").Value = "kWh produced/1kg sludge"
$ws.Range("J39").Value = "`$/kWh of electricity, average US industrial"
$ws.Range("K32").Value = "Metcalf&Eddy"
$ws.Range("K36").Value = "Metcalf&Eddy"
$ws.Range("K37").Value = "See supp calcs"
$ws.Range("L39").Value = $eiaUrl
$ws.Range("G38").Value = "cost of energy recovered by WWTP"

# ---- Numbers / formulas ----
$ws.Range("F30").Value = 1
$ws.Range("I30").Value = 11

$ws.Range("F31").Formula = "=I32*F30"
$ws.Range("I31").Formula = "=AVERAGE(0.1,0.33)"
$ws.Range("I31").NumberFormat = "0.00%"

$ws.Range("F32").Formula = "=F30-F31"
$ws.Range("H32").NumberFormat = "0%"
$ws.Range("I32").Value = 0.59
$ws.Range("I32").NumberFormat = "0%"

$ws.Range("F33").Formula = "=F32/I31"
$ws.Range("I33").Formula = "=I30/I31"

$ws.Range("F34").Formula = "=F33/10^3*I30"
$ws.Range("I34").Formula = "=I33/10^3"

$ws.Range("F35").Formula = "=F31/I37"
$ws.Range("I35").Formula = "=I34/I32"

$ws.Range("F36").Formula = "=F35*I38"
$ws.Range("I36").Value = 0.62
$ws.Range("I36").NumberFormat = "0.00%"

$ws.Range("F37").Formula = "=F36*10^6*0.0000002777778"
$ws.Range("I37").Value = 1.48

$ws.Range("F38").Formula = "=F37*I39"
$ws.Range("I38").Value = 55

$ws.Range("I39").Value = 0.0696

$ws.Range("K39").Value = 43040
$ws.Range("K39").NumberFormat = "mmm-yy"

# Bold "total" style used elsewhere in the sheet (same style as J12/J13)
$ws.Range("J12").Copy() | Out-Null
$ws.Range("F34").PasteSpecial(-4122) | Out-Null
$ws.Range("G34").PasteSpecial(-4122) | Out-Null
$ws.Range("F38").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---- Hyperlink on L39 ----
$ws.Hyperlinks.Add($ws.Range("L39"), $eiaUrl) | Out-Null
# re-apply the workbook's existing "Hyperlink" cell style (same one used by K14/L15/K16/N15)
$ws.Range("K14").Copy() | Out-Null
$ws.Range("L39").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---- View state ----
$ws.Range("G39").Select() | Out-Null

Write-Host "applied supplemental calcs"
